{"js": "// The author merged three runs (\"...we \", \"collapse\", \" onto the grass...\")\n// into a single run of plain text in two places in the document. We locate\n// each target paragraph by the distinctive text it contains and rewrite its\n// content as one run (Office.js's Paragraph.insertText(\"...\", \"Replace\")\n// replaces the whole paragraph's contents, collapsing multiple runs into a\n// single run while keeping the paragraph's/first run's formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst fixes = [\n  {\n    match: \"Thankfully, Prim is about as unathletic\",\n    full:\n      \"Thankfully, Prim is about as unathletic as I am, and after a while I\\u2019m able to catch up to her. Both of us completely out of breath, we collapse onto the grass surrounding the playground she ran to.\",\n  },\n  {\n    match: \"gasp for breath, wishing\",\n    full: \"I gasp for breath, wishing that I had more stamina.\",\n  },\n];\n\nfor (const fix of fixes) {\n  const target = paragraphs.items.find((p) => p.text.indexOf(fix.match) !== -1);\n  if (target) {\n    target.insertText(fix.full, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The author merged three runs (\"...we \", \"collapse\", \" onto the grass...\")\n# into a single run of plain text in two places in the document (same for\n# \"I \" + \"gasp\" + \" for breath...\"). We use Find/Replace with a search string\n# that spans the full, already-correct text of each paragraph: when the\n# matched range crosses run boundaries, Word coalesces the matched text into\n# a single run (identical visible text, identical run formatting), which is\n# exactly the edit described by the diff.\n\n$d = $word.ActiveDocument\n\nfunction Merge-Runs($doc, $searchText) {\n    $rng = $doc.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $searchText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\n$apos = [char]8217\n\n$text1 = \"Thankfully, Prim is about as unathletic as I am, and after a while I\" + $apos + \"m able to catch up to her. Both of us completely out of breath, we collapse onto the grass surrounding the playground she ran to.\"\n$text2 = \"I gasp for breath, wishing that I had more stamina.\"\n\nMerge-Runs $d $text1\nMerge-Runs $d $text2\n"}
